# Penalty Reward System (unfinished) - remove the now-obsolete weekly/monthly
# data points so the remaining rows shift up to take their place.

$wb = $excel.ActiveWorkbook

# "Weekly Quantity" sheet: drop the two stale weekly rows (old rows 4 & 5),
# which shifts the trailing rows (old 6,7,8) up into rows 4,5,6.
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("A4:B5").EntireRow.Delete()

# "Monthly Trend" sheet: drop the stale monthly row (old row 3), which
# shifts the trailing rows (old 4,5) up into rows 3,4.
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("A3:B3").EntireRow.Delete()
